$wb = $excel.ActiveWorkbook

# --- YDS sheet: append new weeks rush/pass yardage-per-play logs ---
$ws = $wb.Worksheets.Item("YDS")
$ws.Range("B2").Value = $ws.Range("B2").Text + " 8 3 3 52 1 6 3 14 3 3 1 4 5 0 14 9 4 12 -1 3 20 7 5 10 0 11 2 0 3 1 38 3 3 2"
$ws.Range("B3").Value = $ws.Range("B3").Text + " 4 19 2 4 4 6 4 22 2 7 12 5"
$ws.Range("C2").Value = $ws.Range("C2").Text + " 3 0 9 3 1 1 6 6 6 0 8 0 2 4 3 2 4 7 1 6 3 9 7 2 3 6 1 0 26"
$ws.Range("C3").Value = $ws.Range("C3").Text + " 1 4 12 9 9 8 17 10 11 1 7 5 10 25 19 11 13 22 18 17 8 7 6 6 4"

# --- OFF sheet: season totals updated through Week 17 ---
$ws = $wb.Worksheets.Item("OFF")
$ws.Range("C2").Value = 360
$ws.Range("D2").Value = 17
$ws.Range("F2").Value = 99
$ws.Range("G2").Value = 120
$ws.Range("J2").Value = 45
$ws.Range("L2").Value = 543
$ws.Range("M2").Value = 322
$ws.Range("O2").Value = 44
$ws.Range("P2").Value = 20
$ws.Range("Q2").Value = 983
$ws.Range("C3").Value = 337
$ws.Range("E3").Value = 60
$ws.Range("F3").Value = 200
$ws.Range("G3").Value = 65
$ws.Range("H3").Value = 45
$ws.Range("I3").Value = 134
$ws.Range("J3").Value = 112
$ws.Range("N3").Value = 42

# --- DEF sheet: season totals updated through Week 17 ---
$ws = $wb.Worksheets.Item("DEF")
$ws.Range("C2").Value = 438
$ws.Range("D2").Value = 30
$ws.Range("F2").Value = 127
$ws.Range("G2").Value = 115
$ws.Range("J2").Value = 61
$ws.Range("L2").Value = 569
$ws.Range("M2").Value = 390
$ws.Range("O2").Value = 31
$ws.Range("P2").Value = 18
$ws.Range("Q2").Value = 1071
$ws.Range("B3").Value = 17
$ws.Range("C3").Value = 331
$ws.Range("D3").Value = 11
$ws.Range("E3").Value = 74
$ws.Range("F3").Value = 205
$ws.Range("G3").Value = 71
$ws.Range("H3").Value = 48
$ws.Range("I3").Value = 114
$ws.Range("J3").Value = 103
$ws.Range("N3").Value = 37

# --- ST sheet: special-teams totals + distance logs ---
$ws = $wb.Worksheets.Item("ST")
$ws.Range("B2").Value = 122
$ws.Range("D2").Value = 136
$ws.Range("F2").Value = 22
$ws.Range("G2").Value = 18
$ws.Range("J2").Value = 15
$ws.Range("K2").Value = 15
$ws.Range("L2").Value = 7
$ws.Range("M2").Value = 4
$ws.Range("B3").Value = 74
$ws.Range("B4").Value = $ws.Range("B4").Text + " 64"
$ws.Range("B5").Value = $ws.Range("B5").Text + " 23"
$ws.Range("B6").Value = $ws.Range("B6").Text + " 20 2 34 31 29 19"
$ws.Range("D3").Value = $ws.Range("D3").Text + " 55"
$ws.Range("D4").Value = $ws.Range("D4").Text + " 4"
$ws.Range("D5").Value = $ws.Range("D5").Text + " 12"

# --- TURNS sheet: turnovers ---
$ws = $wb.Worksheets.Item("TURNS")
$ws.Range("E2").Value = 24
$ws.Range("D3").Value = 6
$ws.Range("E3").Value = 22

# --- PEN sheet: penalties ---
$ws = $wb.Worksheets.Item("PEN")
$ws.Range("B2").Value = 24
$ws.Range("D2").Value = 21
$ws.Range("D4").Value = 19
